# Replace the frequency/temperature range table with new measurement data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(20,22,24,26,28,30,32,34,36,38,40,43,46,49,52,55,58,61,64,67,70,73,76,79,82,85,88,91,94,97,100,103,106,109,112,115,118,121,124,127,130,133,136,139,142,145,148,151,158,165,172,179,186,193,200,207,214,221,228,235,242,249,256,263,270,277,284,291,298,305,312,319,326,333,340,347,354,361,368,375)
$colB = @(50,63,79,100,126,158,200,251,316,398,500,631,794,1000,1259,1585,1995,2512,3162,3981,5012,6310,7943,10000,12590,15850,19950,25120,31620,39810,50120,63100,79430,100000,125900,158500,199500,251200,316200,398100,501200,631000,794300,1000000,1259000,1585000,1995000,2512000,3162000,3981000,5000000)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
}

for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}

# Sort column B (rows 1-51) ascending, matching the sortState recorded in the sheet
$rangeB = $ws.Range("B1:B51")
$null = $ws.Sort.SortFields.Clear()
$null = $ws.Sort.SortFields.Add($rangeB, 0, 1, $null, 0)
$null = $ws.Sort.SetRange($rangeB)
$ws.Sort.Header = 2
$null = $ws.Sort.Apply()

$null = $ws.Range("B1:B51").Select()
Write-Output "done"
